$wb = $excel.ActiveWorkbook

# --- Sheet2 / Table2: add calculated column "Ratio (Contract to Perm)" ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$lo = $ws2.ListObjects.Item(1)
$col = $lo.ListColumns.Add()
$ws2.Range("D1").Value = "Ratio (Contract to Perm)"

$formula = "=Table2[[#This Row],[Contract Employee]]/Table2[[#This Row],[Perm Employee]]"
$ws2.Range("D2").Formula = $formula
$ws2.Range("D3").Formula = $formula
$ws2.Range("D4").Formula = $formula
$ws2.Range("D5").Formula = $formula

$ws2.Columns.Item(2).ColumnWidth = 23.88671875
$ws2.Columns.Item(3).ColumnWidth = 13.33203125

[void]$ws2.Range("D3").Select()

# --- Sheet1: add summary column H "Ratio (Perm to Contract)" ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("H11").Value = "Ratio (Perm to Contract)"

$ws1.Range("H12").FormulaArray = "=SUM(FILTER(Table2[Ratio (Contract to Perm)],Table2[Category]=A12))"
$ws1.Range("H13").FormulaArray = "=SUM(FILTER(Table2[Ratio (Contract to Perm)],Table2[Category]=A13))"
$ws1.Range("H14").FormulaArray = "=SUM(FILTER(Table2[Ratio (Contract to Perm)],Table2[Category]=A14))"
$ws1.Range("H15").FormulaArray = "=SUM(FILTER(Table2[Ratio (Contract to Perm)],Table2[Category]=A15))"

[void]$ws1.Range("E13").Select()
